$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 5854
$ws.Range("E2").Value = 147
$ws.Range("F2").Value = 161
$ws.Range("G2").Value = -98
$ws.Range("H2").Value = -111
$ws.Range("I2").Value = -111
$ws.Range("K2").Value = 5177
$ws.Range("L2").Value = 4180
$ws.Range("M2").Value = 998
$ws.Range("N2").Value = 998
$ws.Range("P2").Value = 277
$ws.Range("Q2").Value = 172
$ws.Range("R2").Value = -620
$ws.Range("S2").Value = 219
$ws.Range("T2").Value = 206
$ws.Range("U2").Value = -33
$ws.Range("V2").Value = 3032
$ws.Range("W2").Value = 2.52
$ws.Range("X2").Value = -1.89
$ws.Range("Y2").Value = -10.67
$ws.Range("Z2").Value = -2.24
$ws.Range("AA2").Value = 418.98
$ws.Range("AB2").Value = 270.43
$ws.Range("AC2").Value = -200
$ws.Range("AD2").Value = -4.75
$ws.Range("AE2").Value = 1815
$ws.Range("AF2").Value = 0.52
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 55320000

# Row 3
$ws.Range("D3").Value = 7468
$ws.Range("E3").Value = 443
$ws.Range("F3").Value = 443
$ws.Range("G3").Value = 240
$ws.Range("H3").Value = 215
$ws.Range("I3").Value = 216
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5611
$ws.Range("L3").Value = 4315
$ws.Range("M3").Value = 1296
$ws.Range("N3").Value = 1294
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 277
$ws.Range("Q3").Value = 442
$ws.Range("R3").Value = -106
$ws.Range("S3").Value = -122
$ws.Range("T3").Value = 543
$ws.Range("U3").Value = -101
$ws.Range("V3").Value = 3006
$ws.Range("W3").Value = 5.93
$ws.Range("X3").Value = 2.88
$ws.Range("Y3").Value = 18.81
$ws.Range("Z3").Value = 3.99
$ws.Range("AA3").Value = 332.92
$ws.Range("AB3").Value = 348.02
$ws.Range("AC3").Value = 390
$ws.Range("AD3").Value = 19.4
$ws.Range("AE3").Value = 2353
$ws.Range("AF3").Value = 3.21
$ws.Range("AG3").Value = 25
$ws.Range("AH3").Value = 0.33
$ws.Range("AI3").Value = 6.38
$ws.Range("AJ3").Value = 55320000

# Row 4
$ws.Range("D4").Value = 10112
$ws.Range("E4").Value = 781
$ws.Range("F4").Value = 781
$ws.Range("G4").Value = 612
$ws.Range("H4").Value = 487
$ws.Range("I4").Value = 445
$ws.Range("J4").Value = 43
$ws.Range("K4").Value = 8039
$ws.Range("L4").Value = 5162
$ws.Range("M4").Value = 2877
$ws.Range("N4").Value = 2307
$ws.Range("O4").Value = 570
$ws.Range("P4").Value = 277
$ws.Range("Q4").Value = 299
$ws.Range("R4").Value = -933
$ws.Range("S4").Value = 1547
$ws.Range("T4").Value = 1063
$ws.Range("U4").Value = -763
$ws.Range("V4").Value = 3519
$ws.Range("W4").Value = 7.73
$ws.Range("X4").Value = 4.82
$ws.Range("Y4").Value = 24.71
$ws.Range("Z4").Value = 7.14
$ws.Range("AA4").Value = 179.46
$ws.Range("AB4").Value = 499.37
$ws.Range("AC4").Value = 804
$ws.Range("AD4").Value = 13
$ws.Range("AE4").Value = 4196
$ws.Range("AF4").Value = 2.49
$ws.Range("AG4").Value = 25
$ws.Range("AH4").Value = 0.24
$ws.Range("AI4").Value = 3.09
$ws.Range("AJ4").Value = 55320000

# Row 5
$ws.Range("D5").Value = 11519
$ws.Range("E5").Value = 948
$ws.Range("F5").Value = 948
$ws.Range("G5").Value = 717
$ws.Range("H5").Value = 485
$ws.Range("I5").Value = 365
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 8757
$ws.Range("L5").Value = 5634
$ws.Range("M5").Value = 3123
$ws.Range("N5").Value = 2423
$ws.Range("O5").Value = 701
$ws.Range("P5").Value = 277
$ws.Range("Q5").Value = 826
$ws.Range("R5").Value = -1860
$ws.Range("S5").Value = 682
$ws.Range("T5").Value = 1645
$ws.Range("U5").Value = -819
$ws.Range("V5").Value = 3883
$ws.Range("W5").Value = 8.23
$ws.Range("X5").Value = 4.21
$ws.Range("Y5").Value = 15.42
$ws.Range("Z5").Value = 5.77
$ws.Range("AA5").Value = 180.39
$ws.Range("AB5").Value = 610.91
$ws.Range("AC5").Value = 659
$ws.Range("AD5").Value = 15.55
$ws.Range("AE5").Value = 4407
$ws.Range("AF5").Value = 2.33
$ws.Range("AG5").Value = 25
$ws.Range("AH5").Value = 0.24
$ws.Range("AI5").Value = 3.77
$ws.Range("AJ5").Value = 55320000

# Row 6
$ws.Range("D6").Value = 11119
$ws.Range("E6").Value = 577
$ws.Range("F6").Value = 577
$ws.Range("G6").Value = 345
$ws.Range("H6").Value = 259
$ws.Range("I6").Value = 203
$ws.Range("K6").Value = 9913
$ws.Range("L6").Value = 6371
$ws.Range("M6").Value = 3541
$ws.Range("N6").Value = 2694
$ws.Range("P6").Value = 277
$ws.Range("Q6").Value = 569
$ws.Range("R6").Value = -848
$ws.Range("S6").Value = 335
$ws.Range("T6").Value = 1160
$ws.Range("U6").Value = -591
$ws.Range("V6").Value = 4511
$ws.Range("W6").Value = 5.19
$ws.Range("X6").Value = 2.33
$ws.Range("Y6").Value = 7.94
$ws.Range("Z6").Value = 2.78
$ws.Range("AA6").Value = 179.92
$ws.Range("AB6").Value = 754.1
$ws.Range("AC6").Value = 367
$ws.Range("AD6").Value = 18.63
$ws.Range("AE6").Value = 4975
$ws.Range("AF6").Value = 1.37
$ws.Range("AG6").Value = 40
$ws.Range("AH6").Value = 0.58
$ws.Range("AI6").Value = 10.66
$ws.Range("AJ6").Value = 55320000

# Row 7
$ws.Range("D7").Value = 13808
$ws.Range("E7").Value = 1099
$ws.Range("G7").Value = 970
$ws.Range("H7").Value = 752
$ws.Range("I7").Value = 576
$ws.Range("K7").Value = 10797
$ws.Range("L7").Value = 6522
$ws.Range("M7").Value = 4274
$ws.Range("N7").Value = 3347
$ws.Range("P7").Value = 277
$ws.Range("Q7").Value = 1221
$ws.Range("R7").Value = -118
$ws.Range("S7").Value = -116
$ws.Range("T7").Value = 300
$ws.Range("W7").Value = 7.96
$ws.Range("X7").Value = 5.44
$ws.Range("Y7").Value = 19.09
$ws.Range("Z7").Value = 7.26
$ws.Range("AA7").Value = 152.59
$ws.Range("AC7").Value = 1042
$ws.Range("AD7").Value = 10.84
$ws.Range("AE7").Value = 6204
$ws.Range("AF7").Value = 1.82
$ws.Range("AG7").Value = 70
$ws.Range("AH7").Value = 0.62
$ws.Range("AI7").Value = 6.72

# Row 8
$ws.Range("D8").Value = 15882
$ws.Range("E8").Value = 1436
$ws.Range("G8").Value = 1296
$ws.Range("H8").Value = 990
$ws.Range("I8").Value = 758
$ws.Range("K8").Value = 11706
$ws.Range("L8").Value = 6578
$ws.Range("M8").Value = 5128
$ws.Range("N8").Value = 4089
$ws.Range("P8").Value = 277
$ws.Range("Q8").Value = 770
$ws.Range("R8").Value = -18
$ws.Range("S8").Value = -158
$ws.Range("T8").Value = 200
$ws.Range("W8").Value = 9.039999999999999
$ws.Range("X8").Value = 6.24
$ws.Range("Y8").Value = 20.4
$ws.Range("Z8").Value = 8.800000000000001
$ws.Range("AA8").Value = 128.27
$ws.Range("AC8").Value = 1371
$ws.Range("AD8").Value = 8.24
$ws.Range("AE8").Value = 7579
$ws.Range("AF8").Value = 1.49
$ws.Range("AG8").Value = 95
$ws.Range("AH8").Value = 0.84
$ws.Range("AI8").Value = 6.93

# Row 9
$ws.Range("D9").Value = 18378
$ws.Range("E9").Value = 1870
$ws.Range("G9").Value = 1742
$ws.Range("H9").Value = 1322
$ws.Range("I9").Value = 1014
$ws.Range("K9").Value = 12953
$ws.Range("L9").Value = 6682
$ws.Range("M9").Value = 6271
$ws.Range("N9").Value = 5072
$ws.Range("P9").Value = 277
$ws.Range("Q9").Value = 1330
$ws.Range("R9").Value = -17
$ws.Range("S9").Value = -216
$ws.Range("T9").Value = 200
$ws.Range("W9").Value = 10.18
$ws.Range("X9").Value = 7.19
$ws.Range("Y9").Value = 22.13
$ws.Range("Z9").Value = 10.72
$ws.Range("AA9").Value = 106.55
$ws.Range("AC9").Value = 1832
$ws.Range("AD9").Value = 6.17
$ws.Range("AE9").Value = 9401
$ws.Range("AF9").Value = 1.2
$ws.Range("AG9").Value = 120
$ws.Range("AH9").Value = 1.06
$ws.Range("AI9").Value = 6.55

# Remove cells that no longer exist in the updated data
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("U9").ClearContents()
